# paises.xlsx update: refresh COVID-19 stats and re-rank several countries
# (Bosnia y Herzegovina, Maldivas, Suazilandia, Sudan del Sur move up in the ranking)
# as of "25 de Abril de 2020 a las 13:22".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Pakistan'
$arr[0,1] = 12188
$arr[0,2] = 248
$arr[0,3] = 2755
$arr[0,4] = 9177
$arr[0,5] = 111
$arr[0,6] = 3
$arr[0,7] = 256
$ws.Range("A32:H32").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Kuwait'
$arr[0,1] = 2892
$arr[0,2] = 278
$arr[0,3] = 656
$arr[0,4] = 2217
$arr[0,5] = 58
$arr[0,6] = 4
$arr[0,7] = 19
$ws.Range("A61:H61").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Bosnia y Herzegovina'
$arr[0,1] = 1486
$arr[0,2] = 65
$arr[0,3] = 592
$arr[0,4] = 837
$arr[0,5] = 4
$arr[0,6] = 2
$arr[0,7] = 57
$ws.Range("A74:H74").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Afganistan'
$arr[0,1] = 1463
$arr[0,2] = 112
$arr[0,3] = 188
$arr[0,4] = 1228
$arr[0,5] = 7
$arr[0,6] = 4
$arr[0,7] = 47
$ws.Range("A75:H75").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Nueva Zelanda'
$arr[0,1] = 1461
$arr[0,2] = 5
$arr[0,3] = 1118
$arr[0,4] = 325
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = 18
$ws.Range("A76:H76").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Camerun'
$arr[0,1] = 1430
$arr[0,2] = 0
$arr[0,3] = 668
$arr[0,4] = 719
$arr[0,5] = 20
$arr[0,6] = 0
$arr[0,7] = 43
$ws.Range("A77:H77").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Lituania'
$arr[0,1] = 1426
$arr[0,2] = 16
$arr[0,3] = 460
$arr[0,4] = 925
$arr[0,5] = 17
$arr[0,6] = 1
$arr[0,7] = 41
$ws.Range("A78:H78").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Maldivas'
$arr[0,1] = 137
$arr[0,2] = 8
$arr[0,3] = 16
$arr[0,4] = 121
$arr[0,5] = 2
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A138:H138").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Gibraltar'
$arr[0,1] = 133
$arr[0,2] = 0
$arr[0,3] = 131
$arr[0,4] = 2
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A139:H139").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Etiopia'
$arr[0,1] = 122
$arr[0,2] = 5
$arr[0,3] = 29
$arr[0,4] = 90
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 3
$ws.Range("A140:H140").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Madagascar'
$arr[0,1] = 122
$arr[0,2] = 0
$arr[0,3] = 61
$arr[0,4] = 61
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A141:H141").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Camboya'
$arr[0,1] = 122
$arr[0,2] = 0
$arr[0,3] = 117
$arr[0,4] = 5
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A142:H142").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Suazilandia'
$arr[0,1] = 40
$arr[0,2] = 4
$arr[0,3] = 10
$arr[0,4] = 29
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 1
$ws.Range("A170:H170").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Puerto Rico'
$arr[0,1] = 39
$arr[0,2] = 0
$arr[0,3] = 1
$arr[0,4] = 36
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 2
$ws.Range("A171:H171").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Eritrea'
$arr[0,1] = 39
$arr[0,2] = 0
$arr[0,3] = 11
$arr[0,4] = 28
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A172:H172").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'San Martin (Parte Francesa)'
$arr[0,1] = 38
$arr[0,2] = 0
$arr[0,3] = 24
$arr[0,4] = 11
$arr[0,5] = 3
$arr[0,6] = 0
$arr[0,7] = 3
$ws.Range("A173:H173").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Mongolia'
$arr[0,1] = 37
$arr[0,2] = 1
$arr[0,3] = 9
$arr[0,4] = 28
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A174:H174").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Sudan del Sur'
$arr[0,1] = 5
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 5
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A210:H210").Value = $arr

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 'Bonaire, San Eustaquio y Saba'
$arr[0,1] = 5
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 5
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$ws.Range("A211:H211").Value = $arr

# Update the "last updated" timestamp shown at the top of the sheet
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 13:22"
